# New Submission Synced: 2026-02-08 19:38:53
# Sheet "JSS 3E" gets its previous last row's Admission No fixed up to a
# real number, and a brand-new form submission is appended as row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3E")

# Row 8 (RAHAMA MOHAMMED YUSUF) had its Admission No stored as text;
# normalize it to a numeric value.
$ws.Range("C8").Value = 42

# Append the newly synced submission as row 9.
$ws.Range("A9").Value = "2026-02-08 19:38:53"
$ws.Range("B9").Value = "GLORIA JOHN GADZAMA"

# Admission No for this submission is kept as text (matches how the sync
# tool originally wrote every other row before they got normalized).
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "35"

$ws.Range("D9").Value = 9
